$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 229
$ws.Range("I2").Value = 161.25
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 161.25
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -48.25
$ws.Range("N2").Value = -726

$ws.Range("H4").Value = 3806.6667
$ws.Range("I4").Value = 2268
$ws.Range("J4").Value = 11500
$ws.Range("K4").Value = 2268
$ws.Range("L4").Value = 11500
$ws.Range("M4").Value = -2154
$ws.Range("N4").Value = -11728

$ws.Range("H5").Value = 236
$ws.Range("I5").Value = 236
$ws.Range("K5").Value = 236
$ws.Range("M5").Value = -121

$ws.Range("H9").Value = 79
$ws.Range("I9").Value = 79
$ws.Range("K9").Value = 79
$ws.Range("M9").Value = 90

$ws.Range("H40").Value = 7061.75
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7061.75
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7061.75
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -7411.75

$ws.Range("H41").Value = 1739.2222
$ws.Range("I41").Value = 3417
$ws.Range("J41").Value = 900.3333
$ws.Range("K41").Value = 3417
$ws.Range("L41").Value = 900.3333
$ws.Range("M41").Value = -2977
$ws.Range("N41").Value = -1780.3333

$ws.Range("H129").Value = 5000
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 2748.6667
$ws.Range("I132").Value = 2748.6667
$ws.Range("K132").Value = 8246.000100000001
$ws.Range("M132").Value = -5716.000100000001

$ws.Range("H137").Value = 2124.75
$ws.Range("I137").Value = 2000
$ws.Range("K137").Value = 6000
$ws.Range("M137").Value = -3450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 582.5
$ws.Range("I4").Value = 582.5
$ws.Range("K4").Value = 582.5
$ws.Range("M4").Value = -466.5

$ws.Range("H32").Value = 1316.4
$ws.Range("I32").Value = 1275
$ws.Range("K32").Value = 1275
$ws.Range("M32").Value = -988

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 384
$ws.Range("I107").Value = 384
$ws.Range("K107").Value = 384
$ws.Range("M107").Value = 1536

$ws.Range("H134").Value = 5008.6665
$ws.Range("I134").Value = 5506
$ws.Range("J134").Value = 4014
$ws.Range("K134").Value = 16518
$ws.Range("L134").Value = 12042
$ws.Range("M134").Value = -13983
$ws.Range("N134").Value = -17112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H135").Value = 99995
$ws.Range("J135").Value = 99995
$ws.Range("L135").Value = 99995
$ws.Range("N135").Value = -110135

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 353263.25
$ws.Range("I4").Value = 333665.06
$ws.Range("K4").Value = 1000995.18
$ws.Range("M4").Value = -1000883.18

$ws.Range("H11").Value = 8335809
$ws.Range("I11").Value = 10002936
$ws.Range("K11").Value = 30008808
$ws.Range("M11").Value = -30008668

$ws.Range("H15").Value = 87.5
$ws.Range("I15").Value = 72.5
$ws.Range("J15").Value = 95
$ws.Range("K15").Value = 217.5
$ws.Range("L15").Value = 285
$ws.Range("M15").Value = -77.5
$ws.Range("N15").Value = -565

$ws.Range("H34").Value = 1038.8182
$ws.Range("J34").Value = 1333.3334
$ws.Range("L34").Value = 4000.0002
$ws.Range("N34").Value = -4168.0002

$ws.Range("I38").Value = 17.333334
$ws.Range("J38").Value = 85.71429000000001
$ws.Range("K38").Value = 52.000002
$ws.Range("L38").Value = 257.14287
$ws.Range("M38").Value = 294.999998
$ws.Range("N38").Value = -951.14287

$ws.Range("H39").Value = 4248.857
$ws.Range("J39").Value = 4482.4546
$ws.Range("L39").Value = 13447.3638
$ws.Range("N39").Value = -14035.3638

$ws.Range("H44").Value = 124.5
$ws.Range("I44").Value = 124.5
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 373.5
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 24.5
$ws.Range("N44").ClearContents()

$ws.Range("H55").Value = 3357.2222
$ws.Range("J55").Value = 3526.875
$ws.Range("L55").Value = 10580.625
$ws.Range("N55").Value = -10934.625

$ws.Range("H104").Value = 1013
$ws.Range("I104").Value = 1013
$ws.Range("K104").Value = 3039
$ws.Range("M104").Value = -418

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2750
$ws.Range("I97").Value = 2750
$ws.Range("K97").Value = 2750
$ws.Range("M97").Value = -2254

$ws.Range("H102").Value = 2577.6365
$ws.Range("I102").Value = 2577.6365
$ws.Range("K102").Value = 2577.6365
$ws.Range("M102").Value = -955.6365000000001

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H124").Value = 99995
$ws.Range("J124").Value = 99995
$ws.Range("L124").Value = 99995
$ws.Range("N124").Value = -109815

$ws.Range("H132").Value = 5125
$ws.Range("J132").Value = 4750
$ws.Range("L132").Value = 14250
$ws.Range("N132").Value = -19310

$ws.Range("H135").Value = 74444
$ws.Range("J135").Value = 74444
$ws.Range("L135").Value = 74444
$ws.Range("N135").Value = -84584

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4640.5557
$ws.Range("I40").Value = 4640.5557
$ws.Range("K40").Value = 4640.5557
$ws.Range("M40").Value = -4504.5557

$ws.Range("H46").Value = 3395.625
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3395.625
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3395.625
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3771.625

$ws.Range("H122").Value = 9500
$ws.Range("I122").Value = 9500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 28500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -26050
$ws.Range("N122").ClearContents()

$ws.Range("H130").Value = 99995
$ws.Range("J130").Value = 99995
$ws.Range("L130").Value = 99995
$ws.Range("N130").Value = -110035

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2047.2222
$ws.Range("I107").Value = 1553.125
$ws.Range("K107").Value = 4659.375
$ws.Range("M107").Value = -2739.375

$ws.Range("H112").Value = 5000
$ws.Range("J112").Value = 5000
$ws.Range("L112").Value = 5000
$ws.Range("N112").Value = -7954

$ws.Range("H113").Value = 781.8421
$ws.Range("I113").Value = 879.75
$ws.Range("J113").Value = 614
$ws.Range("K113").Value = 2639.25
$ws.Range("L113").Value = 1842
$ws.Range("M113").Value = -469.25
$ws.Range("N113").Value = -6182

$ws.Range("H125").Value = 50715
$ws.Range("J125").Value = 50715
$ws.Range("L125").Value = 50715
$ws.Range("N125").Value = -60555

$ws.Range("H131").Value = 99995
$ws.Range("J131").Value = 99995
$ws.Range("L131").Value = 99995
$ws.Range("N131").Value = -110075

$ws.Range("H132").Value = 2311.7
$ws.Range("I132").Value = 1702.125
$ws.Range("K132").Value = 5106.375
$ws.Range("M132").Value = -2576.375
